# od313.xlsx edit:
#  1. "page 1": record the number of poles (49) in the "4) No. of Poles" field (X6).
#  2. "page 3"/"page 4"/"page 5"/"page 6": clear out the stale pole-location data
#     (street, pole #, lat/long, GIS id columns B/D/G/H/I/J/K) for a run of
#     rows while leaving the running pole-count column (A) untouched.

$wb = $excel.ActiveWorkbook

# --- 1. No. of Poles -------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("page 1")
$wsSummary.Range("X6").Value = 49

# --- 2. Clear stale pole-location rows -------------------------------------
$cols = @("B", "D", "G", "H", "I", "J", "K")

function Clear-PoleRows($ws, $rows) {
    foreach ($row in $rows) {
        foreach ($col in $cols) {
            $ws.Range($col + $row).Value = ""
        }
    }
}

# "page 3": pole entries 50-55
$wsPage3 = $wb.Worksheets.Item("page 3")
Clear-PoleRows $wsPage3 @(34, 36, 38, 40, 42, 44)

# "page 4": pole entries 56-75
$wsPage4 = $wb.Worksheets.Item("page 4")
Clear-PoleRows $wsPage4 @(6, 8, 10, 12, 14, 16, 18, 20, 22, 24, 26, 28, 30, 32, 34, 36, 38, 40, 42, 44)

# "page 5": pole entries 76-95
$wsPage5 = $wb.Worksheets.Item("page 5")
Clear-PoleRows $wsPage5 @(6, 8, 10, 12, 14, 16, 18, 20, 22, 24, 26, 28, 30, 32, 34, 36, 38, 40, 42, 44)

# "page 6": pole entries 96-97
$wsPage6 = $wb.Worksheets.Item("page 6")
Clear-PoleRows $wsPage6 @(6, 8)
